$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-price text/percentage cells (safe to assign directly as text)
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +4.95%  "
$ws.Range("E6").Value = "  -6.35%  "
$ws.Range("E7").Value = "  -3.28%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -5.14%  "
$ws.Range("E10").Value = "  -6.73%  "
$ws.Range("E11").Value = "  -8.65%  "
$ws.Range("E12").Value = "  -4.71%  "
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("E15").Value = "  +7.91%  "
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("E21").Value = "  -5.30%  "
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("E23").Value = "  -5.61%  "
$ws.Range("E24").Value = "  -3.92%  "
$ws.Range("E25").Value = "  +4.60%  "
$ws.Range("E26").Value = "  -7.11%  "
$ws.Range("E27").Value = "  -6.38%  "
$ws.Range("E28").Value = "  -5.61%  "
$ws.Range("E29").Value = "  -5.31%  "
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("E31").Value = "  -6.59%  "
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("E33").Value = "  +7.18%  "
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  -9.12%  "
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("E37").Value = "  -9.46%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E40").Value = "  +8.40%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("E44").Value = "  -7.75%  "
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("E46").Value = "  -3.82%  "
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("E48").Value = "  +11.13%  "
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("E51").Value = "  -2.96%  "

# Price cells: force text format so numeric-looking strings are not
# auto-converted/rounded by Excel, then restore the default "Normal" style
# so no stray number-format styling is left on the cell.
$priceCells = @("D2","D3","D5","D6","D7","D8","D10","D12","D13","D14","D15","D16","D20","D21","D22","D23","D24","D25","D26","D28","D29","D30","D33","D34","D35","D36","D37","D38","D39","D40","D44","D45","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.445.66"
$ws.Range("D3").Value = "3.833.01"
$ws.Range("D5").Value = "512.07"
$ws.Range("D6").Value = "138.34"
$ws.Range("D7").Value = "0.602"
$ws.Range("D8").Value = "1.00"
$ws.Range("D10").Value = "0.164"
$ws.Range("D12").Value = "41.10"
$ws.Range("D13").Value = "10.17"
$ws.Range("D14").Value = "4.443.73"
$ws.Range("D15").Value = "21.55"
$ws.Range("D16").Value = "3.857.99"
$ws.Range("D20").Value = "68.473.97"
$ws.Range("D21").Value = "413.45"
$ws.Range("D22").Value = "3.40"
$ws.Range("D23").Value = "13.83"
$ws.Range("D24").Value = "85.92"
$ws.Range("D25").Value = "3.90"
$ws.Range("D26").Value = "11.29"
$ws.Range("D28").Value = "35.08"
$ws.Range("D29").Value = "673.11"
$ws.Range("D30").Value = "13.01"
$ws.Range("D33").Value = "65.49"
$ws.Range("D34").Value = "6.18"
$ws.Range("D35").Value = "0.435"
$ws.Range("D36").Value = "39.27"
$ws.Range("D37").Value = "0.0₃0818"
$ws.Range("D38").Value = "0.999"
$ws.Range("D39").Value = "0.146"
$ws.Range("D40").Value = "3.34"
$ws.Range("D44").Value = "2.73"
$ws.Range("D45").Value = "3.35"
$ws.Range("D47").Value = "2.91"
$ws.Range("D48").Value = "2.699.78"
$ws.Range("D49").Value = "143.40"
$ws.Range("D50").Value = "3.22"

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
